$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.046.05'
$ws.Range('E2').Value = '  -0.74%  '

$ws.Range('D3').Value = '2.018.67'
$ws.Range('E3').Value = '  -1.61%  '

$ws.Range('E4').Value = '  +0.22%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '226.64'
$ws.Range('E5').Value = '  -1.91%  '

$ws.Range('E6').Value = '  -2.15%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '54.74'
$ws.Range('E8').Value = '  -3.99%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.381'
$ws.Range('E9').Value = '  -1.06%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0791'
$ws.Range('E10').Value = '  +2.87%  '

$ws.Range('E11').Value = '  -3.00%  '

$ws.Range('D12').Value = '2.318.23'
$ws.Range('E12').Value = '  -1.61%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.27'
$ws.Range('E13').Value = '  -2.87%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '20.29'
$ws.Range('E14').Value = '  -1.42%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.743'
$ws.Range('E15').Value = '  -1.79%  '

$ws.Range('E16').Value = '  -2.21%  '

$ws.Range('D17').Value = '2.028.90'
$ws.Range('E17').Value = '  -0.74%  '

$ws.Range('D18').Value = '37.032.75'
$ws.Range('E18').Value = '  -0.76%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.10'
$ws.Range('E19').Value = '  +1.91%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '68.83'
$ws.Range('E20').Value = '  -1.32%  '

$ws.Range('D21').Value = '0.0₃0821'
$ws.Range('E21').Value = '  -0.21%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '224.16'
$ws.Range('E22').Value = '  -1.26%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.12%  '

$ws.Range('E24').Value = '  +1.89%  '

$ws.Range('E25').Value = '  -5.44%  '

$ws.Range('E26').Value = '  -2.68%  '

$ws.Range('E27').Value = '  -4.05%  '

$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.36'
$ws.Range('E28').Value = '  +1.42%  '

$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.125'
$ws.Range('E29').Value = '  -3.26%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '18.73'
$ws.Range('E30').Value = '  -2.05%  '

$ws.Range('E31').Value = '  -3.28%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.53'
$ws.Range('E32').Value = '  +0.16%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0614'

$ws.Range('E34').Value = '  -3.00%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.35'
$ws.Range('E35').Value = '  -5.45%  '

$ws.Range('E36').Value = '  +2.15%  '

$ws.Range('E37').Value = '  +0.30%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.13'
$ws.Range('E38').Value = '  -4.51%  '

$ws.Range('E39').Value = '  +2.24%  '

$ws.Range('E40').Value = '  -3.55%  '

$ws.Range('D41').Value = '1.479.16'
$ws.Range('E41').Value = '  -0.47%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '95.32'
$ws.Range('E42').Value = '  -3.23%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '16.55'
$ws.Range('E43').Value = '  -0.65%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0924'
$ws.Range('E44').Value = '  -3.14%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.75'
$ws.Range('E45').Value = '  -5.06%  '

$ws.Range('E46').Value = '  -4.46%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '7.25'
$ws.Range('E47').Value = '  +0.12%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.01'
$ws.Range('E48').Value = '  -1.60%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.91'
$ws.Range('E49').Value = '  -0.71%  '

$ws.Range('D50').Value = '2.203.28'
$ws.Range('E50').Value = '  -1.67%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '44.28'
$ws.Range('E51').Value = '  -2.11%  '
